$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de) on row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C) on row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C) on row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the zh-cn / de-de status columns ---
# Target stored width ~13.41 characters; this engine quantizes ColumnWidth
# assignments onto a 1/6-character grid, so 12.5 is the closest input that
# lands on the nearest achievable grid value to the target.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
